# chamcong_t06_2022 (June 2022 timesheet) — "hoan thanh tinh luong"
#
# 1) Sundays (CN = "Chu Nhat") for June 2022 fall on days 5, 12, 19, 26,
#    which live in columns G, N, U, AB (day N -> column index N+2).
#    The sheet only goes out to day 21 (column W) right now, so the
#    Sundays that matter here are day 5 (G), day 12 (N) and day 19 (U).
#    Every employee's cell on those days becomes the literal "CN" marker
#    instead of "x"/"v".
# 2) The timesheet is extended with 5 new day columns: S..W (days 17-21).
#    Attendance defaults to "x" except the couple of pre-existing "v"
#    (leave) marks that show up in the source data, and the new Sunday
#    column U which is always "CN".
# 3) AH (Tong / total) keeps its COUNTIF(...,"x") formula — forcing a
#    recalculation after the writes refreshes the cached values.
# 4) A fresh conditional-formatting rule highlights "CN" cells (red text)
#    across the whole sheet, re-using the same look as the existing "v"
#    rule. The "v" rule itself is left alone.
# 5) The active selection moves back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 22

$sundayCols = @("G", "N", "U")

# Default attendance mark for the 5 newly added day columns (S..W / days 17-21)
$newCols = @("S", "T", "U", "V", "W")

# Per-row exceptions for the newly added columns (row -> column -> mark)
$exceptions = @{
    2 = @{ "S" = "v" }
    8 = @{ "V" = "v" }
}

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {

    # --- Sunday columns that already existed (G = day5, N = day12) ---
    # Row 22 only has data starting at column R, so skip G/N there.
    if ($row -le 21) {
        $ws.Range("G$row").Value2 = "CN"
        $ws.Range("N$row").Value2 = "CN"
    }

    # --- newly added columns S..W ---
    foreach ($col in $newCols) {
        if ($col -eq "U") {
            $mark = "CN"
        } else {
            $mark = "x"
        }
        if ($exceptions.ContainsKey($row) -and $exceptions[$row].ContainsKey($col)) {
            $mark = $exceptions[$row][$col]
        }
        $ws.Range("$col$row").Value2 = $mark
    }
}

# The "U" column (day 19, a Sunday) sits in the already-centered block of
# the sheet (like R) for the employee rows, so center it the same way.
# Row 22's trailing cells (R..W) were already left in the default/general
# style in the source sheet, so U22 is intentionally skipped here.
$ws.Range("U2:U21").HorizontalAlignment = -4108

# Force the "Tong" (AH) COUNTIF formulas to refresh their cached values.
$excel.CalculateFull()

# --- Conditional formatting: add a "CN" rule, mirroring the "v" one ---
$fcs = $ws.Cells.FormatConditions
$oldCnRule = $fcs.Item(2)
$oldCnRule.Delete()

$cnRule = $fcs.Add(1, 3, '="CN"')
$cnRule.ModifyAppliesToRange($ws.Range("A1:XFD1048576"))
$cnRule.Font.Color = 255
$cnRule.Font.Size = 10
$cnRule.Font.Name = "Calibri"
$cnRule.NumberFormat = "General"

# --- Selection back to A1 ---
$ws.Range("A1").Select()
